$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record is inserted at row 18, pushing the existing
# rows 18-90 down to 19-91 (dimension grows from A1:R90 to A1:R91).
$ws.Rows.Item(18).Insert()

# Populate the newly inserted row with this week's data.
$ws.Cells.Item(18, 1).Value = 6
$ws.Cells.Item(18, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(18, 3).Value = "Metropolitana"
$ws.Cells.Item(18, 4).Value = 45114
$ws.Cells.Item(18, 5).Value = 13
$ws.Cells.Item(18, 6).Value = 100112035
$ws.Cells.Item(18, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(18, 8).Value = "Sin especificar"
$ws.Cells.Item(18, 9).Value = "Primera"
$ws.Cells.Item(18, 10).Value = 410
$ws.Cells.Item(18, 11).Value = 17000
$ws.Cells.Item(18, 12).Value = 18000
$ws.Cells.Item(18, 13).Value = 17439
$ws.Cells.Item(18, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(18, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(18, 16).Value = 1163
$ws.Cells.Item(18, 17).Value = 15
$ws.Cells.Item(18, 18).Value = "Hortaliza"
